# Update the "Förändrad" (Changed) date column (C) for rows 2-27
# from 45328 (2024-02-06) to 45330 (2024-02-08), as produced by the
# automatic update of files.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 27; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45328) {
        $cell.Value2 = 45330
    }
}
